$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shift the three year headers back by one year ---
$ws.Range("A1").Value = 2016
$ws.Range("A13").Value = 2017
$ws.Range("A25").Value = 2018

# --- Add the newly-uploaded 2019 monthly figures (columns C:N, rows 26-35) ---
$data = @{
    26 = @(356,279,312,334,255,249,687,855,573,381,215,300)
    27 = @(277,285,302,286,263,182,505,532,465,392,297,288)
    28 = @(25,15,17,10,10,22,42,55,49,42,37,35)
    29 = @(240,228,253,242,284,334,412,435,402,316,267,256)
    30 = @(88,82,75,79,71,77,91,106,124,121,116,114)
    31 = @(478,483,490,485,476,481,592,660,591,545,541,525)
    32 = @(206,161,157,143,132,126,255,276,259,243,212,201)
    33 = @(645,567,560,489,461,427,561,645,651,612,586,551)
    34 = @(131,100,104,81,92,76,197,206,195,188,172,141)
    35 = @(102,93,88,79,82,80,106,169,139,121,127,121)
}
foreach ($row in $data.Keys) {
    $vals = $data[$row]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($row, 3 + $i).Value = $vals[$i]
    }
}

# --- Column G widened slightly to fit the new values ---
$ws.Columns.Item(7).ColumnWidth = 4.2

# --- Restore the last active selection ---
[void]$ws.Range("N37").Select()
